# Updates cryptos list values (Price and Volume(1h) columns), and fixes two
# row orderings where coin name/link/price/volume were swapped (rows 31-32, 45-46)
# and one coin swap-out (row 51: WhiteBITCoin -> RenderToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.414.68'
$ws.Range("E2").Value = '  +1.94%  '

$ws.Range("D3").Value = '2.621.89'
$ws.Range("E3").Value = '  +1.44%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '569.23'
$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").Value = '142.13'
$ws.Range("E6").Value = '  -0.38%  '

$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").Value = '2.622.07'
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("E11").Value = '  +1.15%  '

$ws.Range("D12").Value = '0.365'
$ws.Range("E12").Value = '  +3.24%  '

$ws.Range("E13").Value = '  -7.23%  '

$ws.Range("D14").Value = '3.065.74'
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").Value = '60.428.37'
$ws.Range("E15").Value = '  +2.04%  '

$ws.Range("D16").Value = '23.34'
$ws.Range("E16").Value = '  +1.51%  '

$ws.Range("E17").Value = '  +2.66%  '

$ws.Range("D18").Value = '2.612.37'
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("D19").Value = '11.35'
$ws.Range("E19").Value = '  +9.52%  '

$ws.Range("E20").Value = '  +1.88%  '

$ws.Range("D21").Value = '346.35'

$ws.Range("E22").Value = '  +8.57%  '

$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").Value = '  +13.62%  '

$ws.Range("D25").Value = '63.23'
$ws.Range("E25").Value = '  -1.31%  '

$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("E27").Value = '  -1.75%  '

$ws.Range("D28").Value = '7.69'
$ws.Range("E28").Value = '  +4.52%  '

$ws.Range("E29").Value = '  +1.20%  '

$ws.Range("D30").Value = '1.82'
$ws.Range("E30").Value = '  +9.10%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '6.37'
$ws.Range("E31").Value = '  +3.91%  '

$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").Value = '160.47'
$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("E34").Value = '  +2.51%  '

$ws.Range("D35").Value = '4.23'
$ws.Range("E35").Value = '  +4.82%  '

$ws.Range("E36").Value = '  +10.34%  '

$ws.Range("D37").Value = '1.21'
$ws.Range("E37").Value = '  +4.31%  '

$ws.Range("E38").Value = '  +8.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.70'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("E40").Value = '  +3.61%  '

$ws.Range("D41").Value = '0.851'
$ws.Range("E41").Value = '  -2.16%  '

$ws.Range("D42").Value = '295.49'
$ws.Range("E42").Value = '  +0.72%  '

$ws.Range("D43").Value = '138.99'
$ws.Range("E43").Value = '  +4.72%  '

$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.0982'
$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.607'
$ws.Range("E46").Value = '  +1.68%  '

$ws.Range("D47").Value = '19.71'
$ws.Range("E47").Value = '  +3.51%  '

$ws.Range("D48").Value = '0.0545'
$ws.Range("E48").Value = '  +1.92%  '

$ws.Range("E49").Value = '  +2.92%  '

$ws.Range("D50").Value = '19.87'
$ws.Range("E50").Value = '  +6.62%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '4.82'
$ws.Range("E51").Value = '  +7.30%  '
